$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 94
$ws.Range("H94").Value = 4232.9414
$ws.Range("I94").Value = 4372.5
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 4372.5
$ws.Range("L94").Value = 2000
$ws.Range("M94").Value = -3921.5
$ws.Range("N94").Value = -2902

# Row 99
$ws.Range("H99").Value = 1490.625
$ws.Range("I99").Value = 1490.625
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4471.875
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -2973.875

# Row 100
$ws.Range("H100").Value = 2691.0908
$ws.Range("I100").Value = 3540
$ws.Range("J100").Value = 1983.6666
$ws.Range("K100").Value = 3540
$ws.Range("L100").Value = 1983.6666
$ws.Range("M100").Value = -2999
$ws.Range("N100").Value = -3065.6666

# Row 101
$ws.Range("H101").Value = 818.3333
$ws.Range("I101").Value = 388.7143
$ws.Range("J101").Value = 1194.25
$ws.Range("K101").Value = 1166.1429
$ws.Range("L101").Value = 3582.75
$ws.Range("M101").Value = 455.8571000000002
$ws.Range("N101").Value = -6826.75

# Row 116
$ws.Range("H116").Value = 8210.888999999999
$ws.Range("I116").Value = 4000
$ws.Range("J116").Value = 9414
$ws.Range("K116").Value = 4000
$ws.Range("L116").Value = 9414
$ws.Range("M116").Value = -558
$ws.Range("N116").Value = -16298

# Row 127
$ws.Range("H127").Value = 981.8148
$ws.Range("I127").Value = 531.875
$ws.Range("J127").Value = 1171.2632
$ws.Range("K127").Value = 1595.625
$ws.Range("L127").Value = 3513.7896
$ws.Range("M127").Value = 3364.375
$ws.Range("N127").Value = -13433.7896

# Row 128
$ws.Range("H128").Value = 43932.668
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 43932.668
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 43932.668
$ws.Range("N128").Value = -53892.668

# Row 135
$ws.Range("H135").Value = 36586540
$ws.Range("I135").Value = 38462440
$ws.Range("J135").Value = 33334976
$ws.Range("K135").Value = 346161960
$ws.Range("L135").Value = 300014784
$ws.Range("M135").Value = -346159425
$ws.Range("N135").Value = -300019854

# Row 140
$ws.Range("H140").Value = 19691.111
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 19691.111
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 19691.111
$ws.Range("N140").Value = -30051.111

# Row 141
$ws.Range("H141").Value = 4640.8125
$ws.Range("I141").Value = 2994.9167
$ws.Range("J141").Value = 9578.5
$ws.Range("K141").Value = 8984.750100000001
$ws.Range("L141").Value = 28735.5
$ws.Range("M141").Value = -3804.750100000001
$ws.Range("N141").Value = -39095.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 26249.652
$ws.Range("I32").Value = 28340.781
$ws.Range("J32").Value = 18034.5
$ws.Range("K32").Value = 28340.781
$ws.Range("L32").Value = 18034.5
$ws.Range("M32").Value = -28053.781
$ws.Range("N32").Value = -18608.5

# Row 110
$ws.Range("H110").Value = 1296.7858
$ws.Range("I110").Value = 1296.7858
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1296.7858
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 748.2141999999999

# Row 139
$ws.Range("H139").Value = 35522.273
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 35522.273
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 35522.273
$ws.Range("N139").Value = -45802.273

# Row 140
$ws.Range("H140").Value = 39366
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 39366
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 39366
$ws.Range("N140").Value = -49726

# Row 141
$ws.Range("H141").Value = 45590.31
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 45590.31
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 45590.31
$ws.Range("N141").Value = -55950.31

$ws = $wb.Worksheets.Item("BSM")
# Row 62
$ws.Range("H62").Value = 20000
$ws.Range("I62").Value = 20000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 20000
$ws.Range("M62").Value = -19314
$ws.Range("L62").ClearContents()
$ws.Range("N62").ClearContents()

# Row 65
$ws.Range("H65").Value = 20000
$ws.Range("I65").Value = 20000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 60000
$ws.Range("M65").Value = -56568
$ws.Range("L65").ClearContents()
$ws.Range("N65").ClearContents()

# Row 81
$ws.Range("H81").Value = 8831.429
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 8831.429
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 8831.429
$ws.Range("N81").Value = -10953.429

# Row 84
$ws.Range("H84").Value = 8831.429
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 8831.429
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 26494.287
$ws.Range("N84").Value = -37102.287

# Row 99
$ws.Range("H99").Value = 1595.5555
$ws.Range("I99").Value = 1531.5834
$ws.Range("J99").Value = 2107.3333
$ws.Range("K99").Value = 1531.5834
$ws.Range("L99").Value = 2107.3333
$ws.Range("M99").Value = -33.58339999999998
$ws.Range("N99").Value = -5103.3333

# Row 138
$ws.Range("H138").Value = 25745.186
$ws.Range("I138").Value = 10000
$ws.Range("J138").Value = 26350.77
$ws.Range("K138").Value = 10000
$ws.Range("L138").Value = 26350.77
$ws.Range("M138").Value = -4860
$ws.Range("N138").Value = -36630.77

# Row 140
$ws.Range("H140").Value = 29180.875
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 29180.875
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 29180.875
$ws.Range("N140").Value = -39540.875

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 149887.83
$ws.Range("I31").Value = 1699.6154
$ws.Range("J31").Value = 204929.16
$ws.Range("K31").Value = 1699.6154
$ws.Range("L31").Value = 204929.16
$ws.Range("M31").Value = -1404.6154
$ws.Range("N31").Value = -205519.16

# Row 34
$ws.Range("H34").Value = 149887.83
$ws.Range("I34").Value = 1699.6154
$ws.Range("J34").Value = 204929.16
$ws.Range("K34").Value = 1699.6154
$ws.Range("L34").Value = 204929.16
$ws.Range("M34").Value = -1497.6154
$ws.Range("N34").Value = -205333.16

# Row 107
$ws.Range("H107").Value = 575.6316
$ws.Range("I107").Value = 438.35715
$ws.Range("J107").Value = 960
$ws.Range("K107").Value = 438.35715
$ws.Range("L107").Value = 960
$ws.Range("M107").Value = 1481.64285
$ws.Range("N107").Value = -4800

# Row 138
$ws.Range("H138").Value = 39923.832
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 39923.832
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 39923.832
$ws.Range("N138").Value = -50203.832

# Row 140
$ws.Range("H140").Value = 28411
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 28411
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 28411
$ws.Range("N140").Value = -38771

# Row 141
$ws.Range("H141").Value = 80664.664
$ws.Range("I141").Value = 129000
$ws.Range("J141").Value = 56497
$ws.Range("K141").Value = 129000
$ws.Range("L141").Value = 56497
$ws.Range("M141").Value = -123820
$ws.Range("N141").Value = -66857

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1520.1154
$ws.Range("I97").Value = 1254.5
$ws.Range("J97").Value = 2405.5
$ws.Range("K97").Value = 1254.5
$ws.Range("L97").Value = 2405.5
$ws.Range("M97").Value = -758.5
$ws.Range("N97").Value = -3397.5

# Row 140
$ws.Range("H140").Value = 30714.9
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 30714.9
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 30714.9
$ws.Range("N140").Value = -41074.9

# Row 141
$ws.Range("H141").Value = 44470
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 44470
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 44470
$ws.Range("N141").Value = -54830

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 2449.0715
$ws.Range("I122").Value = 2382.25
$ws.Range("J122").Value = 2850
$ws.Range("K122").Value = 7146.75
$ws.Range("L122").Value = 8550
$ws.Range("M122").Value = -4696.75
$ws.Range("N122").Value = -13450

# Row 138
$ws.Range("H138").Value = 58666
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 58666
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 58666
$ws.Range("N138").Value = -68946

# Row 140
$ws.Range("H140").Value = 37998.332
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 37998.332
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 37998.332
$ws.Range("N140").Value = -48358.332

# Row 141
$ws.Range("H141").Value = 16602.25
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 16602.25
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 16602.25
$ws.Range("N141").Value = -26962.25

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1627.3334
$ws.Range("I96").Value = 838.4286
$ws.Range("J96").Value = 4388.5
$ws.Range("K96").Value = 838.4286
$ws.Range("L96").Value = 4388.5
$ws.Range("M96").Value = 534.5714
$ws.Range("N96").Value = -7134.5

# Row 122
$ws.Range("H122").Value = 9524476
$ws.Range("I122").Value = 14286214
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 42858642
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -42856192
$ws.Range("N122").Value = -7900

# Row 140
$ws.Range("H140").Value = 49828
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 49828
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 49828
$ws.Range("N140").Value = -60188

# Row 141
$ws.Range("H141").Value = 33735.75
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 33735.75
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 33735.75
$ws.Range("N141").Value = -44095.75
